$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A106").Value = 105
$ws.Range("B106").Value = 1
$ws.Range("C106").Value = "2024-06-17 04:16:40"
$ws.Range("D106").Value = 200
$ws.Range("E106").Value = 7

$ws.Range("A107").Value = 106
$ws.Range("B107").Value = 2
$ws.Range("C107").Value = "2024-06-17 04:16:40"
$ws.Range("D107").Value = 200
$ws.Range("E107").Value = 0
